{"js": "const replacements = [\n  [\"78\u00f79=\", \"15\u00f79=\"],\n  [\"97\u00f72=\", \"72\u00f74=\"],\n  [\"26\u00f73=\", \"74\u00f79=\"],\n  [\"85\u00f78=\", \"83\u00f78=\"],\n  [\"11\u00f77=\", \"77\u00f75=\"],\n  [\"76\u00f75=\", \"40\u00f76=\"],\n  [\"10\u00f78=\", \"21\u00f78=\"],\n  [\"89\u00f76=\", \"68\u00f76=\"],\n  [\"26\u00f76=\", \"46\u00f79=\"],\n  [\"75\u00f75=\", \"91\u00f77=\"],\n  [\"15\u00f76=\", \"13\u00f79=\"],\n  [\"27\u00f72=\", \"99\u00f72=\"],\n  [\"43\u00f73=\", \"67\u00f79=\"],\n  [\"25\u00f73=\", \"87\u00f72=\"],\n  [\"55\u00f74=\", \"83\u00f76=\"],\n  [\"82\u00f79=\", \"21\u00f76=\"],\n  [\"99\u00f79=\", \"82\u00f75=\"],\n  [\"14\u00f74=\", \"64\u00f77=\"],\n  [\"82\u00f74=\", \"25\u00f77=\"],\n  [\"74\u00f77=\", \"28\u00f72=\"],\n  [\"45\u00f79=\", \"74\u00f74=\"],\n  [\"93\u00f74=\", \"18\u00f72=\"],\n  [\"73\u00f78=\", \"87\u00f77=\"],\n  [\"51\u00f72=\", \"36\u00f73=\"],\n  [\"80\u00f77=\", \"52\u00f78=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nfor (const table of tables.items) {\n  const rows = table.rows;\n  rows.load(\"items\");\n  await context.sync();\n\n  for (const row of rows.items) {\n    const cells = row.cells;\n    cells.load(\"items\");\n    await context.sync();\n\n    for (const cell of cells.items) {\n      const paragraphs = cell.body.paragraphs;\n      paragraphs.load(\"items/text\");\n      await context.sync();\n\n      for (const paragraph of paragraphs.items) {\n        const text = paragraph.text;\n        for (const [oldText, newText] of replacements) {\n          if (text === oldText) {\n            paragraph.insertText(newText, Word.InsertLocation.replace);\n            break;\n          }\n        }\n      }\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"78\u00f79=\", \"15\u00f79=\"),\n  @(\"97\u00f72=\", \"72\u00f74=\"),\n  @(\"26\u00f73=\", \"74\u00f79=\"),\n  @(\"85\u00f78=\", \"83\u00f78=\"),\n  @(\"11\u00f77=\", \"77\u00f75=\"),\n  @(\"76\u00f75=\", \"40\u00f76=\"),\n  @(\"10\u00f78=\", \"21\u00f78=\"),\n  @(\"89\u00f76=\", \"68\u00f76=\"),\n  @(\"26\u00f76=\", \"46\u00f79=\"),\n  @(\"75\u00f75=\", \"91\u00f77=\"),\n  @(\"15\u00f76=\", \"13\u00f79=\"),\n  @(\"27\u00f72=\", \"99\u00f72=\"),\n  @(\"43\u00f73=\", \"67\u00f79=\"),\n  @(\"25\u00f73=\", \"87\u00f72=\"),\n  @(\"55\u00f74=\", \"83\u00f76=\"),\n  @(\"82\u00f79=\", \"21\u00f76=\"),\n  @(\"99\u00f79=\", \"82\u00f75=\"),\n  @(\"14\u00f74=\", \"64\u00f77=\"),\n  @(\"82\u00f74=\", \"25\u00f77=\"),\n  @(\"74\u00f77=\", \"28\u00f72=\"),\n  @(\"45\u00f79=\", \"74\u00f74=\"),\n  @(\"93\u00f74=\", \"18\u00f72=\"),\n  @(\"73\u00f78=\", \"87\u00f77=\"),\n  @(\"51\u00f72=\", \"36\u00f73=\"),\n  @(\"80\u00f77=\", \"52\u00f78=\")\n)\n\nforeach ($pair in $replacements) {\n  $find = $pair[0]\n  $replace = $pair[1]\n\n  $rng = $d.Content\n  $rng.Find.ClearFormatting()\n  $rng.Find.Replacement.ClearFormatting()\n  $rng.Find.Text = $find\n  $rng.Find.Replacement.Text = $replace\n  $rng.Find.Forward = $true\n  $rng.Find.Wrap = 1\n  $rng.Find.MatchWholeWord = $false\n  $rng.Find.MatchWildcards = $false\n  $rng.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)\n}\n"}
